$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('3153 Sacred Heart Community St Kilda Tier 1B Outbreak', 15),
    @('3175 The Bays Aged Care Facility Hastings Outbreak', 15),
    @('3600 Belvedere Aged Care Noble Park Outbreak', 27),
    @('3612 BlueCross Glengowrie Outbreak', 35),
    @('3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak', 31),
    @('3980 Arcare Keysborough Aged Care Keysborough Outbreak', 19),
    @('4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak', 17),
    @('4518 Regis Aged Care Fawkner Outbreak', 13),
    @('AW Window Transport Group Depot North Geelong Outbreak', 10),
    @('Allied Pinnacle Factory Altona North Outbreak', 15),
    @('Bespoke Childcare Dingley Village Outbreak', 10),
    @('Bread Solutions Braeside Outbreak', 15),
    @('CS Square Caroline Springs Outbreak', 11),
    @('Child''s Play Early Learning Centre Tarneit Outbreak', 11),
    @('Community Kids Pascoe Vale Early Education Centre Pascoe Vale Outbreak', 22),
    @('Essential Caravans Somerton Outbreak', 10),
    @('Guardian Childcare Caulfield Outbreak', 21),
    @('Hello Fresh Warehouse Ravenhall Outbreak', 17),
    @('Inghams Enterprises Somerville Outbreak', 24),
    @('Lantmannen Unibake Australia Mordialloc Outbreak', 21),
    @('Launch Housing City Edge Crisis Accommodation South Melbourne Outbreak', 11),
    @('MacKillop Family Services Residential Facility Glenroy Outbreak', 10),
    @('Melbourne Assessment Prison West Melbourne Outbreak', 16),
    @('Melbourne Youth Justice Centre Parkville Outbreak', 10),
    @('Mildura Police Station Mildura Outbreak', 10),
    @('Monash Health Kingston Centre South 5 Outbreak', 10),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B Outbreak', 43),
    @('Northern Health The Northern Hospital Epping Outbreak', 19),
    @('Oceania Meat Processors Laverton North Outbreak', 17),
    @('Public Housing 140 Brunswick Street Fitzroy Outbreak', 11),
    @('Shawlands Caravan Park Dandenong South Outbreak', 17),
    @('St Vincents Hospital Emergency Department Melbourne Outbreak', 37),
    @('Target Distribution Centre Truganina Outbreak', 11),
    @('The Robin Hood Inn Drouin West Outbreak', 48),
    @('The Toolshed Bar Private Event Noojee Outbreak', 16),
    @('Turosi Breakwater Outbreak', 10),
    @('Visy Recycling Springvale Outbreak', 14),
    @('Werribee Mercy Hospital Emergency Department Outbreak', 26),
    @('Western Health Footscray Hospital Emergency Department Outbreak', 11),
    @('Western Health Sunshine Hospital Emergency Department Outbreak', 24)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
